$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price (column D) and volume-change (column E) figures.
# Ensure the target cells keep their original text formatting (e.g. "30.582.18")
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.582.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.129.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5283"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4498"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.06"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09395"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.185"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.42"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.720"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.004"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.101.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "103.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001174"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.63"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06733"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.371"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.006"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.594.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.79"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.339"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.364.04"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.550"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.168"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.789"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +10.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1066"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +13.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.321"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.968"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.71"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02670"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06902"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7148"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.77"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2267"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.336"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6978"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.82"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.401"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.289"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +10.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.640"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.237"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.68%  "
